$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 65, shifting existing rows 65-68 down to 66-69
$ws.Rows.Item(65).Insert()

# Fill in the new row 65 with data
$ws.Cells.Item(65, 1).Value = 10
$ws.Cells.Item(65, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(65, 3).Value = "La Araucanía"
$ws.Cells.Item(65, 4).Value = 45005
$ws.Cells.Item(65, 5).Value = 9
$ws.Cells.Item(65, 6).Value = 100112010
$ws.Cells.Item(65, 7).Value = "Achicoria"
$ws.Cells.Item(65, 8).Value = "Sin especificar"
$ws.Cells.Item(65, 9).Value = "Primera"
$ws.Cells.Item(65, 10).Value = 150
$ws.Cells.Item(65, 11).Value = 10000
$ws.Cells.Item(65, 12).Value = 10000
$ws.Cells.Item(65, 13).Value = 10000
$ws.Cells.Item(65, 14).Value = "$/caja 18 unidades"
$ws.Cells.Item(65, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(65, 16).Value = 556
$ws.Cells.Item(65, 17).Value = 18
$ws.Cells.Item(65, 18).Value = "Hortaliza"
